$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rng, $val) {
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '68.621.85'
$ws.Range("E2").Value = '  +0.91%  '

Set-TextValue $ws.Range("D3") '3.696.20'
$ws.Range("E3").Value = '  -3.00%  '

Set-TextValue $ws.Range("D4") '0.999'
$ws.Range("E4").Value = '  -0.07%  '

Set-TextValue $ws.Range("D5") '599.56'
$ws.Range("E5").Value = '  +0.41%  '

Set-TextValue $ws.Range("D6") '166.84'
$ws.Range("E6").Value = '  -4.66%  '

Set-TextValue $ws.Range("D7") '3.694.46'
$ws.Range("E7").Value = '  -3.06%  '

$ws.Range("E8").Value = '  -0.03%  '

Set-TextValue $ws.Range("D9") '0.533'
$ws.Range("E9").Value = '  +1.01%  '

$ws.Range("E10").Value = '  +2.62%  '

Set-TextValue $ws.Range("D11") '6.29'
$ws.Range("E11").Value = '  -0.08%  '

Set-TextValue $ws.Range("D12") '0.458'
$ws.Range("E12").Value = '  -1.55%  '

Set-TextValue $ws.Range("D13") '38.00'
$ws.Range("E13").Value = '  -0.07%  '

$ws.Range("E14").Value = '  -1.08%  '

Set-TextValue $ws.Range("D15") '4.314.75'
$ws.Range("E15").Value = '  -2.91%  '

Set-TextValue $ws.Range("D16") '3.694.04'
$ws.Range("E16").Value = '  -2.97%  '

Set-TextValue $ws.Range("D17") '68.594.65'
$ws.Range("E17").Value = '  +0.77%  '

Set-TextValue $ws.Range("D18") '7.24'
$ws.Range("E18").Value = '  +1.23%  '

$ws.Range("E19").Value = '  -0.63%  '

Set-TextValue $ws.Range("D20") '17.08'
$ws.Range("E20").Value = '  +4.85%  '

Set-TextValue $ws.Range("D21") '491.54'
$ws.Range("E21").Value = '  +0.32%  '

Set-TextValue $ws.Range("D22") '9.11'
$ws.Range("E22").Value = '  -1.19%  '

Set-TextValue $ws.Range("D23") '0.720'
$ws.Range("E23").Value = '  -1.79%  '

Set-TextValue $ws.Range("D24") '84.41'
$ws.Range("E24").Value = '  -0.43%  '

Set-TextValue $ws.Range("D25") '0.0000141'
$ws.Range("E25").Value = '  +1.23%  '

Set-TextValue $ws.Range("D26") '2.29'
$ws.Range("E26").Value = '  -4.39%  '

Set-TextValue $ws.Range("D27") '12.18'
$ws.Range("E27").Value = '  -1.10%  '

Set-TextValue $ws.Range("D28") '10.05'
$ws.Range("E28").Value = '  -1.81%  '

$ws.Range("E29").Value = '  +0.06%  '

Set-TextValue $ws.Range("D30") '2.92'
$ws.Range("E30").Value = '  -0.22%  '

$ws.Range("E31").Value = '  +0.25%  '

Set-TextValue $ws.Range("D32") '2.37'
$ws.Range("E32").Value = '  -3.30%  '

Set-TextValue $ws.Range("D33") '31.39'
$ws.Range("E33").Value = '  -4.17%  '

Set-TextValue $ws.Range("D34") '3.832.44'
$ws.Range("E34").Value = '  -2.99%  '

$ws.Range("E35").Value = '  -0.82%  '

Set-TextValue $ws.Range("D36") '3.636.80'
$ws.Range("E36").Value = '  -3.01%  '

Set-TextValue $ws.Range("D37") '0.999'
$ws.Range("E37").Value = '  -0.06%  '

Set-TextValue $ws.Range("D38") '0.999'
$ws.Range("E38").Value = '  -0.98%  '

Set-TextValue $ws.Range("D39") '5.72'
$ws.Range("E39").Value = '  -1.23%  '

$ws.Range("E40").Value = '  -4.03%  '

Set-TextValue $ws.Range("D41") '0.321'
$ws.Range("E41").Value = '  -1.65%  '

Set-TextValue $ws.Range("D42") '49.08'
$ws.Range("E42").Value = '  +0.32%  '

Set-TextValue $ws.Range("D43") '431.84'
$ws.Range("E43").Value = '  -4.02%  '

Set-TextValue $ws.Range("D44") '1.96'
$ws.Range("E44").Value = '  -1.76%  '

Set-TextValue $ws.Range("D45") '2.81'
$ws.Range("E45").Value = '  -3.16%  '

Set-TextValue $ws.Range("D46") '8.37'
$ws.Range("E46").Value = '  +0.95%  '

$ws.Range("E47").Value = '  +0.00%  '

Set-TextValue $ws.Range("D48") '40.17'
$ws.Range("E48").Value = '  -3.43%  '

Set-TextValue $ws.Range("D49") '141.49'
$ws.Range("E49").Value = '  +1.64%  '

Set-TextValue $ws.Range("D50") '0.0349'
$ws.Range("E50").Value = '  -0.85%  '

Set-TextValue $ws.Range("D51") '2.731.60'
$ws.Range("E51").Value = '  -3.60%  '
